# -----------------------------------------------------------------------
# Updated symbol list on Sun Jan 29 14:34:34 UTC 2023 with GitHub Actions
#
# This applies refreshed Price / Volume(1h) figures (and, for rows 8-17,
# the coin/link that rotated into that row) to the "cryptos" worksheet.
# All target cells are plain text (e.g. "314.61", "2.84%"), so the number
# format is forced to Text ("@") before writing the value - otherwise Excel
# would silently reinterpret them as a number / percentage.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "314.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.84%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.23%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.115"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.53%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08186"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.58%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.045"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.29%"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.258"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.92%"

# Row 8
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9324"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.34%"

# Row 9
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1414"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.73%"

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1991"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.71%"

# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09117"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.93%"

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03520"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.17%"

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09807"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.19%"

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001398"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.52%"

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006322"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.58%"

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.656"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.77%"

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.275"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.18%"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3461"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.03%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1305"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.12%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.902"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.16%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2451"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04330"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.46%"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001225"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.50%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004777"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "15.99%"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.08%"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004001"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-10.03%"

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.73%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05250"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.46%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007546"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.96%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009777"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.31%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1376"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.10%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002126"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.76%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009468"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.17%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006462"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.51%"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.07%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002768"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.49%"

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-25.05%"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.07%"

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.07%"
